# Apply the price-list corrections for rows 17-25 (brand/price swaps) and
# update the SalesQuantity (K) / Turnover (L) rollups for rows 32-35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 ---
$ws.Range("A17").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("G17").Value = 3.98
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 30
$ws.Range("J17").Value = "Pantene"
$ws.Range("L17").Value = 9

# --- Row 18 ---
$ws.Range("G18").Value = 7.95
$ws.Range("J18").Value = "Dixan"
$ws.Range("L18").Value = 17.96

# --- Row 19 ---
$ws.Range("A19").Value = "Πελάτες Τιμή Πώλησης"
$ws.Range("G19").Value = 1.75
$ws.Range("H19").Value = 0.99
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = "Farmer"
$ws.Range("K19").Value = 5
$ws.Range("L19").Value = 4.4

# --- Row 20 ---
$ws.Range("G20").Value = 1.55
$ws.Range("H20").Value = 1.15
$ws.Range("J20").Value = "Bazaar"
$ws.Range("L20").Value = 4.65

# --- Row 23 ---
$ws.Range("A23").Value = "Πωλήσεις Έκπτωση 1"
$ws.Range("G23").Value = 2.99
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = "Colgate"
$ws.Range("K23").Value = 7
$ws.Range("L23").Value = 8.34

# --- Row 24 ---
$ws.Range("L24").Value = 8.460000000000001

# --- Row 25 ---
$ws.Range("A25").Value = "Πελάτες Τιμή Πώλησης"
$ws.Range("G25").Value = 2.1
$ws.Range("H25").Value = 1.69
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = "Farmer"
$ws.Range("K25").Value = 8
$ws.Range("L25").Value = 13.8

# --- Row 32 ---
$ws.Range("K32").Value = 20
$ws.Range("L32").Value = 52.2

# --- Row 33 ---
$ws.Range("K33").Value = 49
$ws.Range("L33").Value = 56.17

# --- Row 34 ---
$ws.Range("K34").Value = 114
$ws.Range("L34").Value = 54.96

# --- Row 35 (totals) ---
$ws.Range("K35").Value = 337
$ws.Range("L35").Value = 413.5599999999999
